$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new row 6: batch_005 hash-partition dml case.
# (Typed in the same left-to-right column order the author used, so any
# newly introduced shared strings land in a sensible sequence.) ---
$ws.Range("A6").Value = "batch_005"
$ws.Range("B6").Value = "y"
$ws.Range("C6").Value = "批量操作语句5执行"
$ws.Range("D6").Value = "batchsql"
$ws.Range("F6").Value = "batch05"
$ws.Range("H6").Value = "batch_sql_05"
$ws.Range("I6").Value = "select * from `$batch05"
$ws.Range("M6").Value = "csv_containsAll"

# --- Fix row 5 (batch_004): Query_result1 path was pointing at the wrong
# folder ("cases/batchsql" instead of "mysqlcases/batchsql"); correct it. ---
$ws.Range("J5").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_004.csv"

# --- Finish row 6 ---
$ws.Range("J6").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_005.csv"

# --- Add new row 7: batch_006 hash-partition dml case ---
$ws.Range("A7").Value = "batch_006"
$ws.Range("B7").Value = "y"
$ws.Range("C7").Value = "批量操作语句6执行"
$ws.Range("D7").Value = "batchsql"
$ws.Range("F7").Value = "batch06"
$ws.Range("H7").Value = "batch_sql_06"
$ws.Range("I7").Value = "select * from `$batch06"
$ws.Range("J7").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_006.csv"
$ws.Range("M7").Value = "csv_containsAll"

# --- Match the formatting used throughout the sheet: every populated data
# cell is stored as text, and the Query_result columns (J) use "fill"
# horizontal alignment like the rest of column J. ---
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","M")) {
    $ws.Range("${col}6").NumberFormat = "@"
}
foreach ($col in @("A","B","C","D","F","H","I","J","M")) {
    $ws.Range("${col}7").NumberFormat = "@"
}
$ws.Range("J6").HorizontalAlignment = 5
$ws.Range("J7").HorizontalAlignment = 5

# Update the active selection to match the author's final cursor position
$ws.Range("G6").Select() | Out-Null
